# Apply the row-permutation edit described by the diff: each data row (2-16)
# in the "Artfynd" sheet is repointed at a different underlying record, so
# columns A,B,D,E,F,G,H,I,Q,R,S are rewritten cell-by-cell (only where the
# value actually differs from what is already there).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Cells.Item(2, 1).Value = 111902028
$ws.Cells.Item(2, 2).Value = 90666
$ws.Cells.Item(2, 4).NumberFormat = "@"
$ws.Cells.Item(2, 4).Value = 'LC'
$ws.Cells.Item(2, 5).Value = 4364
$ws.Cells.Item(2, 6).NumberFormat = "@"
$ws.Cells.Item(2, 6).Value = 'Dropptaggsvamp'
$ws.Cells.Item(2, 7).NumberFormat = "@"
$ws.Cells.Item(2, 7).Value = 'Hydnellum ferrugineum'
$ws.Cells.Item(2, 8).NumberFormat = "@"
$ws.Cells.Item(2, 8).Value = '(Fr.:Fr.) P. Karst.'
$ws.Cells.Item(2, 9).NumberFormat = "@"
$ws.Cells.Item(2, 9).Value = '1'
$ws.Cells.Item(2, 17).Value = 524954.0254130038
$ws.Cells.Item(2, 18).Value = 6867304.187839299
# Row 3
$ws.Cells.Item(3, 1).Value = 111902038
$ws.Cells.Item(3, 17).Value = 524892.725176702
$ws.Cells.Item(3, 18).Value = 6867498.641564975
$ws.Cells.Item(3, 19).Value = 10
# Row 4
$ws.Cells.Item(4, 1).Value = 111902035
$ws.Cells.Item(4, 2).Value = 90658
$ws.Cells.Item(4, 5).Value = 4361
$ws.Cells.Item(4, 6).NumberFormat = "@"
$ws.Cells.Item(4, 6).Value = 'Orange taggsvamp'
$ws.Cells.Item(4, 7).NumberFormat = "@"
$ws.Cells.Item(4, 7).Value = 'Hydnellum aurantiacum'
$ws.Cells.Item(4, 8).NumberFormat = "@"
$ws.Cells.Item(4, 8).Value = '(Batsch:Fr.) P.Karst.'
$ws.Cells.Item(4, 9).NumberFormat = "@"
$ws.Cells.Item(4, 9).Value = '3'
$ws.Cells.Item(4, 17).Value = 525047.2558985724
$ws.Cells.Item(4, 18).Value = 6867385.376238698
$ws.Cells.Item(4, 19).Value = 25
# Row 5
$ws.Cells.Item(5, 1).Value = 111902036
$ws.Cells.Item(5, 2).Value = 88032
$ws.Cells.Item(5, 4).NumberFormat = "@"
$ws.Cells.Item(5, 4).Value = 'VU'
$ws.Cells.Item(5, 5).Value = 6276
$ws.Cells.Item(5, 6).NumberFormat = "@"
$ws.Cells.Item(5, 6).Value = 'Goliatmusseron'
$ws.Cells.Item(5, 7).NumberFormat = "@"
$ws.Cells.Item(5, 7).Value = 'Tricholoma matsutake'
$ws.Cells.Item(5, 8).NumberFormat = "@"
$ws.Cells.Item(5, 8).Value = '(S.Ito & S.Imai) Singer'
$ws.Cells.Item(5, 17).Value = 525015.987664115
$ws.Cells.Item(5, 18).Value = 6867405.860822954
# Row 6
$ws.Cells.Item(6, 1).Value = 111902037
$ws.Cells.Item(6, 2).Value = 90654
$ws.Cells.Item(6, 4).NumberFormat = "@"
$ws.Cells.Item(6, 4).Value = 'VU'
$ws.Cells.Item(6, 5).Value = 149
$ws.Cells.Item(6, 6).NumberFormat = "@"
$ws.Cells.Item(6, 6).Value = 'Tallgråticka'
$ws.Cells.Item(6, 7).NumberFormat = "@"
$ws.Cells.Item(6, 7).Value = 'Boletopsis grisea'
$ws.Cells.Item(6, 8).NumberFormat = "@"
$ws.Cells.Item(6, 8).Value = '(Peck) Bondartsev & Singer'
$ws.Cells.Item(6, 9).NumberFormat = "@"
$ws.Cells.Item(6, 9).Value = '2'
$ws.Cells.Item(6, 17).Value = 524868.6293626219
$ws.Cells.Item(6, 18).Value = 6867441.031870116
$ws.Cells.Item(6, 19).Value = 5
# Row 7
$ws.Cells.Item(7, 1).Value = 111902033
$ws.Cells.Item(7, 2).Value = 90300
$ws.Cells.Item(7, 5).Value = 4745
$ws.Cells.Item(7, 6).NumberFormat = "@"
$ws.Cells.Item(7, 6).Value = 'Tallriska'
$ws.Cells.Item(7, 7).NumberFormat = "@"
$ws.Cells.Item(7, 7).Value = 'Lactarius musteus'
$ws.Cells.Item(7, 8).NumberFormat = "@"
$ws.Cells.Item(7, 8).Value = 'Fr.'
$ws.Cells.Item(7, 9).NumberFormat = "@"
$ws.Cells.Item(7, 9).Value = '1'
$ws.Cells.Item(7, 17).Value = 525027.0938798942
$ws.Cells.Item(7, 18).Value = 6867370.16309081
$ws.Cells.Item(7, 19).Value = 10
# Row 8
$ws.Cells.Item(8, 1).Value = 111902031
$ws.Cells.Item(8, 2).Value = 90660
$ws.Cells.Item(8, 5).Value = 4362
$ws.Cells.Item(8, 6).NumberFormat = "@"
$ws.Cells.Item(8, 6).Value = 'Blå taggsvamp'
$ws.Cells.Item(8, 7).NumberFormat = "@"
$ws.Cells.Item(8, 7).Value = 'Hydnellum caeruleum'
$ws.Cells.Item(8, 8).NumberFormat = "@"
$ws.Cells.Item(8, 8).Value = '(Hornem.) P.Karst.'
$ws.Cells.Item(8, 9).NumberFormat = "@"
$ws.Cells.Item(8, 9).Value = '2'
$ws.Cells.Item(8, 17).Value = 524990.2026765908
$ws.Cells.Item(8, 18).Value = 6867385.898910107
$ws.Cells.Item(8, 19).Value = 25
# Row 9
$ws.Cells.Item(9, 1).Value = 111902029
$ws.Cells.Item(9, 2).Value = 88032
$ws.Cells.Item(9, 5).Value = 6276
$ws.Cells.Item(9, 6).NumberFormat = "@"
$ws.Cells.Item(9, 6).Value = 'Goliatmusseron'
$ws.Cells.Item(9, 7).NumberFormat = "@"
$ws.Cells.Item(9, 7).Value = 'Tricholoma matsutake'
$ws.Cells.Item(9, 8).NumberFormat = "@"
$ws.Cells.Item(9, 8).Value = '(S.Ito & S.Imai) Singer'
$ws.Cells.Item(9, 9).NumberFormat = "@"
$ws.Cells.Item(9, 9).Value = '4'
$ws.Cells.Item(9, 17).Value = 524971.6686743505
$ws.Cells.Item(9, 18).Value = 6867341.509407703
# Row 10
$ws.Cells.Item(10, 1).Value = 111902027
$ws.Cells.Item(10, 2).Value = 90660
$ws.Cells.Item(10, 5).Value = 4362
$ws.Cells.Item(10, 6).NumberFormat = "@"
$ws.Cells.Item(10, 6).Value = 'Blå taggsvamp'
$ws.Cells.Item(10, 7).NumberFormat = "@"
$ws.Cells.Item(10, 7).Value = 'Hydnellum caeruleum'
$ws.Cells.Item(10, 8).NumberFormat = "@"
$ws.Cells.Item(10, 8).Value = '(Hornem.) P.Karst.'
$ws.Cells.Item(10, 9).NumberFormat = "@"
$ws.Cells.Item(10, 9).Value = '5'
$ws.Cells.Item(10, 17).Value = 524936.9216418237
$ws.Cells.Item(10, 18).Value = 6867321.952660743
$ws.Cells.Item(10, 19).Value = 25
# Row 11
$ws.Cells.Item(11, 1).Value = 111902030
$ws.Cells.Item(11, 2).Value = 88032
$ws.Cells.Item(11, 4).NumberFormat = "@"
$ws.Cells.Item(11, 4).Value = 'VU'
$ws.Cells.Item(11, 5).Value = 6276
$ws.Cells.Item(11, 6).NumberFormat = "@"
$ws.Cells.Item(11, 6).Value = 'Goliatmusseron'
$ws.Cells.Item(11, 7).NumberFormat = "@"
$ws.Cells.Item(11, 7).Value = 'Tricholoma matsutake'
$ws.Cells.Item(11, 8).NumberFormat = "@"
$ws.Cells.Item(11, 8).Value = '(S.Ito & S.Imai) Singer'
$ws.Cells.Item(11, 9).NumberFormat = "@"
$ws.Cells.Item(11, 9).Value = '6'
$ws.Cells.Item(11, 17).Value = 524971.3961406752
$ws.Cells.Item(11, 18).Value = 6867378.699329315
$ws.Cells.Item(11, 19).Value = 5
# Row 12
$ws.Cells.Item(12, 1).Value = 111902032
$ws.Cells.Item(12, 2).Value = 90658
$ws.Cells.Item(12, 4).NumberFormat = "@"
$ws.Cells.Item(12, 4).Value = 'NT'
$ws.Cells.Item(12, 5).Value = 4361
$ws.Cells.Item(12, 6).NumberFormat = "@"
$ws.Cells.Item(12, 6).Value = 'Orange taggsvamp'
$ws.Cells.Item(12, 7).NumberFormat = "@"
$ws.Cells.Item(12, 7).Value = 'Hydnellum aurantiacum'
$ws.Cells.Item(12, 8).NumberFormat = "@"
$ws.Cells.Item(12, 8).Value = '(Batsch:Fr.) P.Karst.'
$ws.Cells.Item(12, 9).NumberFormat = "@"
$ws.Cells.Item(12, 9).Value = '1'
$ws.Cells.Item(12, 17).Value = 524989.2701192262
$ws.Cells.Item(12, 18).Value = 6867384.479730026
$ws.Cells.Item(12, 19).Value = 5
# Row 13
$ws.Cells.Item(13, 1).Value = 111902026
$ws.Cells.Item(13, 2).Value = 90682
$ws.Cells.Item(13, 4).NumberFormat = "@"
$ws.Cells.Item(13, 4).Value = 'NT'
$ws.Cells.Item(13, 5).Value = 2059
$ws.Cells.Item(13, 6).NumberFormat = "@"
$ws.Cells.Item(13, 6).Value = 'Skrovlig taggsvamp'
$ws.Cells.Item(13, 7).NumberFormat = "@"
$ws.Cells.Item(13, 7).Value = 'Hydnellum scabrosum'
$ws.Cells.Item(13, 8).NumberFormat = "@"
$ws.Cells.Item(13, 8).Value = '(Fr.) E.Larss., K.H.Larss. & Kõljalg'
$ws.Cells.Item(13, 17).Value = 524951.0483835863
$ws.Cells.Item(13, 18).Value = 6867324.410012136
# Row 14
$ws.Cells.Item(14, 1).Value = 111902040
$ws.Cells.Item(14, 2).Value = 90300
$ws.Cells.Item(14, 4).NumberFormat = "@"
$ws.Cells.Item(14, 4).Value = 'NT'
$ws.Cells.Item(14, 5).Value = 4745
$ws.Cells.Item(14, 6).NumberFormat = "@"
$ws.Cells.Item(14, 6).Value = 'Tallriska'
$ws.Cells.Item(14, 7).NumberFormat = "@"
$ws.Cells.Item(14, 7).Value = 'Lactarius musteus'
$ws.Cells.Item(14, 8).NumberFormat = "@"
$ws.Cells.Item(14, 8).Value = 'Fr.'
$ws.Cells.Item(14, 9).NumberFormat = "@"
$ws.Cells.Item(14, 9).Value = '1'
$ws.Cells.Item(14, 17).Value = 524890.9316995766
$ws.Cells.Item(14, 18).Value = 6866840.436305572
$ws.Cells.Item(14, 19).Value = 10
# Row 15
$ws.Cells.Item(15, 1).Value = 111902039
$ws.Cells.Item(15, 2).Value = 90682
$ws.Cells.Item(15, 5).Value = 2059
$ws.Cells.Item(15, 6).NumberFormat = "@"
$ws.Cells.Item(15, 6).Value = 'Skrovlig taggsvamp'
$ws.Cells.Item(15, 7).NumberFormat = "@"
$ws.Cells.Item(15, 7).Value = 'Hydnellum scabrosum'
$ws.Cells.Item(15, 8).NumberFormat = "@"
$ws.Cells.Item(15, 8).Value = '(Fr.) E.Larss., K.H.Larss. & Kõljalg'
$ws.Cells.Item(15, 9).NumberFormat = "@"
$ws.Cells.Item(15, 9).Value = '5'
$ws.Cells.Item(15, 17).Value = 524868.0170565489
$ws.Cells.Item(15, 18).Value = 6867460.329015278
$ws.Cells.Item(15, 19).Value = 5
# Row 16
$ws.Cells.Item(16, 1).Value = 111902034
$ws.Cells.Item(16, 2).Value = 90660
$ws.Cells.Item(16, 5).Value = 4362
$ws.Cells.Item(16, 6).NumberFormat = "@"
$ws.Cells.Item(16, 6).Value = 'Blå taggsvamp'
$ws.Cells.Item(16, 7).NumberFormat = "@"
$ws.Cells.Item(16, 7).Value = 'Hydnellum caeruleum'
$ws.Cells.Item(16, 8).NumberFormat = "@"
$ws.Cells.Item(16, 8).Value = '(Hornem.) P.Karst.'
$ws.Cells.Item(16, 9).NumberFormat = "@"
$ws.Cells.Item(16, 9).Value = '10'
$ws.Cells.Item(16, 17).Value = 525038.6070930503
$ws.Cells.Item(16, 18).Value = 6867407.439287313
